# Add a new "2022-Q4" quarterly sheet right after "总计", shifting the
# existing quarter sheets (2022-Q3 .. 2020-Q4) one slot to the right, and
# record the new quarter's summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item(1)          # "总计"
$templateSheet = $wb.Worksheets.Item(2)         # currently "2022-Q3" - same column layout we need

# --- 1. Create the new "2022-Q4" sheet by duplicating the 2022-Q3 layout ---
$templateSheet.Copy($null, $summarySheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The template (2022-Q3) has 13 data rows (rows 2-14); our new sheet only
# needs 6, so drop the extra rows and shift the remaining cells up.
$newSheet.Range("A8:H14").Delete(-4162)

# --- 2. Fill in the 2022-Q4 fund holdings data ---
$q4Data = @(
    @("000965", "汇丰晋信新动力混合", "0.95", "91.04", "4.69", "0.0446", 3),
    @("540004", "汇丰晋信2026周期混合", "1.08", "23.97", "1.90", "0.0205", 4),
    @("008082", "国寿安保研究精选混合A", "0.34", "84.56", "3.62", "0.0123", 7),
    @("007316", "交银施罗德可转债债券A", "0.82", "20.54", "0.68", "0.0056", 5),
    @("008083", "国寿安保研究精选混合C", "0.13", "84.56", "3.62", "0.0047", 7),
    @("007317", "交银施罗德可转债债券C", "0.22", "20.54", "0.68", "0.0015", 5)
)

for ($i = 0; $i -lt $q4Data.Count; $i++) {
    $row = $i + 2
    $vals = $q4Data[$i]

    $newSheet.Range("A$row").Value = $i

    $newSheet.Range("B$row").Value = "'" + $vals[0]
    $newSheet.Range("B$row").Style = "Normal"

    $newSheet.Range("C$row").Value = $vals[1]
    $newSheet.Range("C$row").Style = "Normal"

    $newSheet.Range("D$row").Value = "'" + $vals[2]
    $newSheet.Range("D$row").Style = "Normal"

    $newSheet.Range("E$row").Value = "'" + $vals[3]
    $newSheet.Range("E$row").Style = "Normal"

    $newSheet.Range("F$row").Value = "'" + $vals[4]
    $newSheet.Range("F$row").Style = "Normal"

    $newSheet.Range("G$row").Value = "'" + $vals[5]
    $newSheet.Range("G$row").Style = "Normal"

    $newSheet.Range("H$row").Value = $vals[6]
    $newSheet.Range("H$row").Style = "Normal"
}

# --- 3. Insert the 2022-Q4 summary row at the top of the "总计" sheet ---
$summarySheet.Rows.Item(2).Insert()

# Re-apply the index-column style (bordered/centered) to the new A2 cell,
# copying the formatting from the row below.
$summarySheet.Range("A3").Copy()
$summarySheet.Range("A2").PasteSpecial(-4122)

# Now (re)write every data row so the full table - including the rows that
# got shifted down by the insert - ends up with the correct values.
$totals = @(
    @("2022-Q4", 6, 0.09),
    @("2022-Q3", 13, 1.14),
    @("2022-Q2", 10, 3.88),
    @("2022-Q1", 27, 25.37),
    @("2021-Q4", 19, 24.79),
    @("2021-Q3", 30, 26.43),
    @("2021-Q2", 33, 18.7),
    @("2021-Q1", 34, 20.18),
    @("2020-Q4", 6, 11.03)
)

for ($i = 0; $i -lt $totals.Count; $i++) {
    $row = $i + 2
    $vals = $totals[$i]

    $summarySheet.Range("A$row").Value = $i

    $summarySheet.Range("B$row").Value = $vals[0]
    $summarySheet.Range("B$row").Style = "Normal"

    $summarySheet.Range("C$row").Value = $vals[1]
    $summarySheet.Range("C$row").Style = "Normal"

    $summarySheet.Range("D$row").Value = $vals[2]
    $summarySheet.Range("D$row").Style = "Normal"
}

# Restore the "总计" sheet as the active tab (matches the original view state).
$summarySheet.Select()

